$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.417.86"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.379.63"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'575.36"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'137.38"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.378.07"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'7.48"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "3.952.85"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'0.124"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "'26.18"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").Value = "3.374.20"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "61.491.79"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'14.04"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'5.86"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'9.37"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "'376.75"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").Value = "3.508.37"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +7.65%  "
$ws.Range("D27").Value = "'71.61"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "'1.73"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "'8.29"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "'23.64"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "'5.29"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("D37").Value = "'6.84"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'165.75"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").Value = "'0.0776"
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  +5.22%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.776"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'41.56"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.41"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'24.44"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").Value = "'6.84"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "'22.74"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "2.363.60"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").Value = "'2.39"
$ws.Range("E51").Value = "  -1.66%  "
